$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H87").Value = 49083.2
$ws.Range("J87").Value = 49083.2
$ws.Range("L87").Value = 49083.2
$ws.Range("N87").Value = -51579.2
$ws.Range("H90").Value = 49083.2
$ws.Range("J90").Value = 49083.2
$ws.Range("L90").Value = 147249.6
$ws.Range("N90").Value = -159729.6
$ws.Range("H98").Value = 1546.5454
$ws.Range("I98").Value = 1535.3334
$ws.Range("J98").Value = 1560
$ws.Range("K98").Value = 1535.3334
$ws.Range("L98").Value = 1560
$ws.Range("M98").Value = -37.33339999999998
$ws.Range("N98").Value = -4556
$ws.Range("H107").Value = 10973
$ws.Range("I107").Value = 10973
$ws.Range("K107").Value = 10973
$ws.Range("M107").Value = -9053
$ws.Range("H111").Value = 1601.8636
$ws.Range("I111").Value = 2201.125
$ws.Range("J111").Value = 1259.4286
$ws.Range("K111").Value = 6603.375
$ws.Range("L111").Value = 3778.2858
$ws.Range("M111").Value = -3536.375
$ws.Range("N111").Value = -9912.2858
$ws.Range("H112").Value = 1376.4
$ws.Range("J112").Value = 1387.9166
$ws.Range("L112").Value = 4163.7498
$ws.Range("N112").Value = -6379.7498
$ws.Range("H113").Value = 243201.83
$ws.Range("I113").Value = 479077
$ws.Range("J113").Value = 7326.6665
$ws.Range("K113").Value = 479077
$ws.Range("L113").Value = 7326.6665
$ws.Range("M113").Value = -475823
$ws.Range("N113").Value = -13834.6665
$ws.Range("H115").Value = 1495
$ws.Range("I115").Value = 742.5
$ws.Range("J115").Value = 3000
$ws.Range("K115").Value = 2227.5
$ws.Range("L115").Value = 9000
$ws.Range("M115").Value = -660.5
$ws.Range("N115").Value = -12134
$ws.Range("H116").Value = 6349.4
$ws.Range("I116").Value = 5174.8335
$ws.Range("J116").Value = 8111.25
$ws.Range("K116").Value = 5174.8335
$ws.Range("L116").Value = 8111.25
$ws.Range("M116").Value = -1732.8335
$ws.Range("N116").Value = -14995.25
$ws.Range("H118").Value = 2459
$ws.Range("I118").Value = 1289.9333
$ws.Range("J118").Value = 4212.6
$ws.Range("K118").Value = 3869.7999
$ws.Range("L118").Value = 12637.8
$ws.Range("M118").Value = -2212.7999
$ws.Range("N118").Value = -15951.8
$ws.Range("H122").Value = 1546.5454
$ws.Range("I122").Value = 1535.3334
$ws.Range("J122").Value = 1560
$ws.Range("K122").Value = 4606.0002
$ws.Range("L122").Value = 4680
$ws.Range("M122").Value = -2156.0002
$ws.Range("N122").Value = -9580
$ws.Range("H137").Value = 5558194
$ws.Range("I137").Value = 1635.591
$ws.Range("K137").Value = 4906.772999999999
$ws.Range("M137").Value = -2356.772999999999
$ws.Range("H138").Value = 5955413
$ws.Range("J138").Value = 14711226
$ws.Range("L138").Value = 44133678
$ws.Range("N138").Value = -44143958

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 3408.7144
$ws.Range("I45").Value = 3572.2
$ws.Range("K45").Value = 3572.2
$ws.Range("M45").Value = -3195.2

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 9263946
$ws.Range("I31").Value = 8332.177
$ws.Range("K31").Value = 8332.177
$ws.Range("M31").Value = -8037.177
$ws.Range("H34").Value = 9263946
$ws.Range("I34").Value = 8332.177
$ws.Range("K34").Value = 8332.177
$ws.Range("M34").Value = -8130.177

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 892.8
$ws.Range("I107").Value = 337.54544
$ws.Range("J107").Value = 1571.4445
$ws.Range("K107").Value = 1012.63632
$ws.Range("L107").Value = 4714.333500000001
$ws.Range("M107").Value = 907.36368
$ws.Range("N107").Value = -8554.333500000001
$ws.Range("H109").Value = 2544.8462
$ws.Range("I109").Value = 1196.909
$ws.Range("J109").Value = 3533.3333
$ws.Range("K109").Value = 3590.727
$ws.Range("L109").Value = 10599.9999
$ws.Range("M109").Value = -2550.727
$ws.Range("N109").Value = -12679.9999
$ws.Range("H112").Value = 5783.778
$ws.Range("J112").Value = 6666.6665
$ws.Range("L112").Value = 19999.9995
$ws.Range("N112").Value = -22215.9995
$ws.Range("H115").Value = 2987.8
$ws.Range("I115").Value = 1319.5
$ws.Range("J115").Value = 4100
$ws.Range("K115").Value = 3958.5
$ws.Range("L115").Value = 12300
$ws.Range("M115").Value = -2783.5
$ws.Range("N115").Value = -14650
$ws.Range("H118").Value = 1558.1538
$ws.Range("I118").Value = 537.4
$ws.Range("K118").Value = 1612.2
$ws.Range("M118").Value = -369.1999999999998
$ws.Range("H131").Value = 814.55
$ws.Range("I131").Value = 466.3889
$ws.Range("J131").Value = 890.9756
$ws.Range("K131").Value = 1399.1667
$ws.Range("L131").Value = 2672.9268
$ws.Range("M131").Value = 3640.8333
$ws.Range("N131").Value = -12752.9268

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 10000
$ws.Range("J95").Value = 10000
$ws.Range("L95").Value = 10000
$ws.Range("N95").Value = -15492
$ws.Range("H113").Value = 2298.8333
$ws.Range("I113").Value = 1945
$ws.Range("K113").Value = 1945
$ws.Range("M113").Value = 225
$ws.Range("H132").Value = 6480.2104
$ws.Range("I132").Value = 5841.5
$ws.Range("J132").Value = 7575.143
$ws.Range("K132").Value = 17524.5
$ws.Range("L132").Value = 22725.429
$ws.Range("M132").Value = -14994.5
$ws.Range("N132").Value = -27785.429

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 898.75
$ws.Range("I46").Value = 758
$ws.Range("J46").Value = 1133.3334
$ws.Range("K46").Value = 758
$ws.Range("L46").Value = 1133.3334
$ws.Range("M46").Value = -570
$ws.Range("N46").Value = -1509.3334
$ws.Range("H55").Value = 482.33334
$ws.Range("I55").Value = 250
$ws.Range("J55").Value = 598.5
$ws.Range("K55").Value = 250
$ws.Range("L55").Value = 598.5
$ws.Range("M55").Value = -77
$ws.Range("N55").Value = -944.5
$ws.Range("H58").Value = 7600
$ws.Range("I58").Value = 7150
$ws.Range("J58").Value = 8500
$ws.Range("K58").Value = 7150
$ws.Range("L58").Value = 8500
$ws.Range("M58").Value = -6890
$ws.Range("N58").Value = -9020
$ws.Range("H130").Value = 43552.668
$ws.Range("J130").Value = 43552.668
$ws.Range("L130").Value = 43552.668
$ws.Range("N130").Value = -53592.668
